$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metrics data for rows 2-26 (columns A-Q)
# Column order: A(model name), B..Q metric values
$data = @(
    @("model_40_5_24", 0.9999686932010947, 0.9990683788130419, 0.9998878524160755, 0.9996339721860903, 0.9997939221234302, 0.00002922352256840288, 0.0008696274845815399, 0.00006235954508296577, 0.0001195570870914823, 0.00009095878724233114, 0.0003955008296977327, 0.005405878519575045, 1.000022768581022, 0.005636017959864586, 134.8810732111485, 204.3569952286359),
    @("model_40_5_23", 0.9999687973865142, 0.9990681690414399, 0.9998884550606907, 0.9996353409956713, 0.9997948048387553, 0.00002912627005247983, 0.000869823297164095, 0.000062024445183913, 0.0001191099874993943, 0.00009056917377781438, 0.0003974780177185774, 0.005396875953038001, 1.000022692809808, 0.005626632135432179, 134.8877400801541, 204.3636620976416),
    @("model_40_5_22", 0.9999689108880594, 0.9990679221086182, 0.9998891529852887, 0.9996368191541218, 0.9997957958334236, 0.00002902032134223802, 0.0008700537981140858, 0.00006163636495150548, 0.0001186271708611983, 0.00009013176790635186, 0.0003996211171519534, 0.005387051265974552, 1.000022610263229, 0.005616389191100532, 134.8950284748787, 204.3709504923661),
    @("model_40_5_21", 0.9999690225315013, 0.9990676288582035, 0.9998898827729913, 0.999638368964277, 0.9997968292274875, 0.00002891610708974671, 0.0008703275345038698, 0.00006123056727356664, 0.0001181209503482621, 0.00008967564776208876, 0.0004020066474911803, 0.005377369904492968, 1.000022529067999, 0.005606295674017497, 134.9022235587118, 204.3781455761992),
    @("model_40_5_20", 0.9999691871756303, 0.9990673112750342, 0.9998909091083912, 0.9996405637970817, 0.9997982875802472, 0.00002876241902246547, 0.0008706239844521262, 0.00006065987456316815, 0.0001174040435810729, 0.00008903195907212054, 0.0004043707607809485, 0.005363060602162301, 1.000022409326814, 0.005591377195062274, 134.9128818359851, 204.3888038534726),
    @("model_40_5_19", 0.9999693181839072, 0.9990669091093074, 0.9998918473370044, 0.9996425302382493, 0.9997996061783132, 0.00002864012854651893, 0.0008709993885051826, 0.00006013817353793935, 0.0001167617372616877, 0.00008844995539981351, 0.0004071815233235464, 0.005351647274112797, 1.000022314048068, 0.005579477978008856, 134.9214034550516, 204.397325472539),
    @("model_40_5_18", 0.999969522832272, 0.9990664643367928, 0.9998932172803393, 0.9996453606655741, 0.9998015171538295, 0.00002844909828099117, 0.0008714145641242795, 0.00005937641799973748, 0.0001158372237867217, 0.00008760648778308106, 0.0004098305136356565, 0.005333769612665246, 1.000022165212893, 0.005560839227501622, 134.9347881855255, 204.4107102030129),
    @("model_40_5_17", 0.9999697673920718, 0.9990660006289223, 0.9998948346415368, 0.9996488561735012, 0.9998038285568004, 0.00002822081244290091, 0.0008718474150669663, 0.00005847708602145658, 0.0001146954724503809, 0.00008658627923591874, 0.0004126017507709247, 0.005312326462379822, 1.000021987351221, 0.005538483197915106, 134.9509016446805, 204.4268236621679),
    @("model_40_5_16", 0.9999700424296993, 0.9990654361960747, 0.9998967057535783, 0.999652902057502, 0.9998065059638899, 0.00002796407689029078, 0.0008723742883544885, 0.00005743665615555715, 0.0001133739496385557, 0.00008540452355272004, 0.0004154343506122972, 0.005288107117891125, 1.000021787323855, 0.005513232785790571, 134.9691796794026, 204.44510169689),
    @("model_40_5_15", 0.9999703378015345, 0.9990648067613277, 0.9998987751084882, 0.9996575437860619, 0.9998095214551393, 0.00002768836024747063, 0.000872961838061817, 0.00005628599355287509, 0.0001118578038031807, 0.00008407354406304106, 0.0004185163132064266, 0.005261973037508899, 1.000021572507975, 0.00548598612350147, 134.988996881511, 204.4649188989984),
    @("model_40_5_14", 0.9999706851514797, 0.9990641067380178, 0.9999012152314812, 0.9996631126455016, 0.9998131227647623, 0.0000273641243206045, 0.0008736152790940616, 0.00005492916575095374, 0.0001100388256060794, 0.00008248399567851655, 0.0004214802748357441, 0.00523107296074185, 1.000021319889833, 0.00545377056649447, 135.0125554677166, 204.488477485204),
    @("model_40_5_13", 0.9999711757214924, 0.999063328765914, 0.9999044236509672, 0.9996707919505045, 0.9998179851884746, 0.00002690619874727947, 0.0008743414819040436, 0.00005314512749903131, 0.0001075305043743543, 0.00008033781593669282, 0.0004237966971188752, 0.005187118539929416, 1.000020963111642, 0.005407944914989856, 135.0463077224022, 204.5222297398896),
    @("model_40_5_12", 0.9999717384115622, 0.9990624267183169, 0.9999081067572059, 0.9996799453957844, 0.999823692068362, 0.00002638095226634647, 0.0008751835037405758, 0.00005109714018177619, 0.0001045406789760495, 0.00007781890957891285, 0.0004257205142285146, 0.005136239116936289, 1.0000205538825, 0.005354899449624987, 135.0857366265694, 204.5616586440568),
    @("model_40_5_11", 0.9999723031702517, 0.9990613886243148, 0.9999119537711741, 0.9996901065481058, 0.9998298750588185, 0.00002585377481975956, 0.0008761525189244472, 0.00004895801214537101, 0.0001012217023112472, 0.00007508985722830911, 0.0004276332186712626, 0.005084660737921416, 1.000020143148908, 0.005301125272233488, 135.1261078796543, 204.6020298971418),
    @("model_40_5_10", 0.9999729218015967, 0.9990603125920591, 0.9999160647267517, 0.9997022770491455, 0.9998369677921808, 0.00002527630961392475, 0.0008771569478027678, 0.00004667211965708867, 0.00009724640426704456, 0.00007195926196206662, 0.0004293971567955426, 0.005027555033405875, 1.000019693235202, 0.005241588459651615, 135.1712859596073, 204.6472079770947),
    @("model_40_5_9", 0.9999735428137912, 0.9990590723539515, 0.9999203813927658, 0.9997159916968761, 0.9998447614576101, 0.00002469662198956527, 0.0008783146556360934, 0.00004427184209873836, 0.00009276673558928505, 0.00006851928884401171, 0.0004307264652256297, 0.004969569598020061, 1.00001924158997, 0.005181134464235028, 135.2176881704782, 204.6936101879656),
    @("model_40_5_8", 0.9999742321481195, 0.9990577015270294, 0.9999250507964877, 0.9997325265397023, 0.9998538195764162, 0.00002405315865241038, 0.0008795942623954167, 0.00004167542511213187, 0.00008736589562931953, 0.00006452121047186993, 0.0004310809287207883, 0.004904401966846761, 1.000018740255913, 0.005113192512087087, 135.2704884662758, 204.7464104837632),
    @("model_40_5_7", 0.9999750739735731, 0.9990562497391828, 0.9999303684323184, 0.9997530572485867, 0.9998647677352113, 0.00002326735154341249, 0.0008809494426241004, 0.00003871855934905888, 0.00008065987041253965, 0.00005968890501960907, 0.0004289500754641852, 0.004823624316156109, 1.00001812801922, 0.005028975989573745, 135.336918802668, 204.8128408201554),
    @("model_40_5_6", 0.9999758841629195, 0.9990546727603843, 0.9999354535167394, 0.9997763235373071, 0.999876578962492, 0.00002251107535179807, 0.0008824214831112251, 0.00003589100355065942, 0.0000730603121244982, 0.00005447565783757881, 0.0004251580252114594, 0.004744583791208463, 1.000017538790604, 0.004946570545842635, 135.403006264011, 204.8789282814985),
    @("model_40_5_5", 0.9999766055490805, 0.9990530095337073, 0.9999401528899602, 0.9998021049521673, 0.9998890771763023, 0.00002183769303567055, 0.0008839740322069808, 0.00003327792205604926, 0.00006463922841272106, 0.00004895918809421269, 0.000419442750750297, 0.004673081749303189, 1.000017014146123, 0.004872024513983822, 135.4637460858034, 204.9396681032908),
    @("model_40_5_4", 0.9999773980277492, 0.9990511089049402, 0.9999449145465017, 0.9998317735515114, 0.9999030557112186, 0.00002109794898416914, 0.000885748185733119, 0.0000306302079869412, 0.00005494845852888681, 0.000042789333257914, 0.0004075515330722631, 0.004593250372467099, 1.000016437798001, 0.004788794550162088, 135.5326694534973, 205.0085914709848),
    @("model_40_5_3", 0.9999776673673897, 0.999048590799919, 0.9999479114752124, 0.9998594221550012, 0.9999151738339821, 0.00002084653225244451, 0.0008880987262383104, 0.00002896376895624547, 0.00004591748773988307, 0.00003744062834806427, 0.0003947577203827783, 0.004565800286088355, 1.000016241914626, 0.00476017585677669, 135.5566458847591, 205.0325679022466),
    @("model_40_5_2", 0.999977532158921, 0.9990452423296278, 0.999950002419727, 0.9998857437883851, 0.9999262298447242, 0.00002097274342300587, 0.0008912243762743194, 0.00002780110147687505, 0.00003731995035258656, 0.00003256071913322549, 0.0003759108334251262, 0.004579600792973757, 1.000016340248058, 0.00477456387981119, 135.5445737910445, 205.020495808532),
    @("model_40_5_1", 0.999976606230205, 0.9990403156318033, 0.9999505290477907, 0.9999069613707922, 0.9999344128028888, 0.00002183705723578099, 0.0008958232324365053, 0.00002750827050069192, 0.00003038956896817542, 0.00002894891973443367, 0.0003538341529279341, 0.004673013720906561, 1.00001701365076, 0.004871953589477512, 135.463804316235, 204.9397263337224),
    @("model_40_5_0", 0.9999746664217333, 0.9990329343000769, 0.9999497592151534, 0.9999210992064256, 0.9999391590820464, 0.00002364778329638102, 0.0009027133816000839, 0.00002793633512209189, 0.00002577167278137618, 0.00002685400395173403, 0.0003294315470545197, 0.004862898651666619, 1.000018424420558, 0.00506992231485602, 135.3044823539183, 204.7804043714057)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowData[$j]
    }
}